# Weekly update: insert the newest week's price records for
# "Terminal Hortofrutícola Agro Chillán - Choclo" ahead of the existing
# history (rows shift down by two to make room for the two new records).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data block (old rows 228-244) down by two rows so the
# two new weekly records can be inserted at the top of the block.
$ws.Rows.Item(228).Insert()
$ws.Rows.Item(228).Insert()

# New record 1: Choclero / Primera
$ws.Range("A228").Value = 7
$ws.Range("B228").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C228").Value = "Ñuble"
$ws.Range("D228").Value = 44918
$ws.Range("E228").Value = 16
$ws.Range("F228").Value = 100112024
$ws.Range("G228").Value = "Choclo"
$ws.Range("H228").Value = "Choclero"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 20000
$ws.Range("K228").Value = 300
$ws.Range("L228").Value = 350
$ws.Range("M228").Value = 325
$ws.Range("N228").Value = "$/unidad"
$ws.Range("O228").Value = "Región de O'Higgins"
$ws.Range("P228").Value = 325
$ws.Range("Q228").Value = 1
$ws.Range("R228").Value = "Hortaliza"

# New record 2: Choclero / Segunda
$ws.Range("A229").Value = 7
$ws.Range("B229").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C229").Value = "Ñuble"
$ws.Range("D229").Value = 44918
$ws.Range("E229").Value = 16
$ws.Range("F229").Value = 100112024
$ws.Range("G229").Value = "Choclo"
$ws.Range("H229").Value = "Choclero"
$ws.Range("I229").Value = "Segunda"
$ws.Range("J229").Value = 15000
$ws.Range("K229").Value = 250
$ws.Range("L229").Value = 250
$ws.Range("M229").Value = 250
$ws.Range("N229").Value = "$/unidad"
$ws.Range("O229").Value = "Región de O'Higgins"
$ws.Range("P229").Value = 250
$ws.Range("Q229").Value = 1
$ws.Range("R229").Value = "Hortaliza"
